$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking "Price" (D) values so Excel
# does not silently reinterpret them as numbers (they are text in the
# source workbook, e.g. thousands-grouped "28.382.17").
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range('D2') '28.382.17'
$ws.Range('E2').Value = '  +3.39%  '

# Row 3
Set-TextValue $ws.Range('D3') '1.863.49'
$ws.Range('E3').Value = '  +1.95%  '

# Row 4
$ws.Range('E4').Value = '  -0.57%  '

# Row 5
Set-TextValue $ws.Range('D5') '337.01'
$ws.Range('E5').Value = '  +1.86%  '

# Row 6
$ws.Range('E6').Value = '  -0.53%  '

# Row 7
Set-TextValue $ws.Range('D7') '0.4698'
$ws.Range('E7').Value = '  +2.59%  '

# Row 8
Set-TextValue $ws.Range('D8') '0.3963'
$ws.Range('E8').Value = '  +3.81%  '

# Row 9
Set-TextValue $ws.Range('D9') '47.58'
$ws.Range('E9').Value = '  +2.40%  '

# Row 10
Set-TextValue $ws.Range('D10') '0.08005'
$ws.Range('E10').Value = '  +1.30%  '

# Row 11
Set-TextValue $ws.Range('D11') '0.9941'
$ws.Range('E11').Value = '  +2.68%  '

# Row 12
$ws.Range('E12').Value = '  +4.13%  '

# Row 13
Set-TextValue $ws.Range('D13') '6.018'
$ws.Range('E13').Value = '  +2.57%  '

# Row 14
Set-TextValue $ws.Range('D14') '1.866.71'
$ws.Range('E14').Value = '  +0.96%  '

# Row 15
Set-TextValue $ws.Range('D15') '7.244'
$ws.Range('E15').Value = '  +2.90%  '

# Row 16
Set-TextValue $ws.Range('D16') '90.18'
$ws.Range('E16').Value = '  +2.21%  '

# Row 17
Set-TextValue $ws.Range('D17') '1.003'
$ws.Range('E17').Value = '  -0.46%  '

# Row 18
Set-TextValue $ws.Range('D18') '0.00001036'
$ws.Range('E18').Value = '  +0.57%  '

# Row 19
Set-TextValue $ws.Range('D19') '0.06611'
$ws.Range('E19').Value = '  -0.60%  '

# Row 20
Set-TextValue $ws.Range('D20') '17.49'
$ws.Range('E20').Value = '  +1.90%  '

# Row 21
$ws.Range('E21').Value = '  -0.44%  '

# Row 22
Set-TextValue $ws.Range('D22') '28.387.30'
$ws.Range('E22').Value = '  +3.49%  '

# Row 23
Set-TextValue $ws.Range('D23') '5.459'
$ws.Range('E23').Value = '  +2.54%  '

# Row 24
Set-TextValue $ws.Range('D24') '11.00'
$ws.Range('E24').Value = '  +2.26%  '

# Row 25
Set-TextValue $ws.Range('D25') '2.271'
$ws.Range('E25').Value = '  -1.38%  '

# Row 26
Set-TextValue $ws.Range('D26') '2.091.19'
$ws.Range('E26').Value = '  +1.50%  '

# Row 27
Set-TextValue $ws.Range('D27') '161.28'
$ws.Range('E27').Value = '  +2.66%  '

# Row 28
$ws.Range('E28').Value = '  +1.73%  '

# Row 29
Set-TextValue $ws.Range('D29') '2.108'
$ws.Range('E29').Value = '  +2.38%  '

# Row 30
Set-TextValue $ws.Range('D30') '5.449'
$ws.Range('E30').Value = '  +4.22%  '

# Row 31
Set-TextValue $ws.Range('D31') '119.21'
$ws.Range('E31').Value = '  +0.83%  '

# Row 32
Set-TextValue $ws.Range('D32') '0.09513'
$ws.Range('E32').Value = '  +2.37%  '

# Row 33
Set-TextValue $ws.Range('D33') '0.9587'
$ws.Range('E33').Value = '  +1.42%  '

# Row 34
Set-TextValue $ws.Range('D34') '3.594'
$ws.Range('E34').Value = '  +0.04%  '

# Row 35
$ws.Range('B35').Value = 'Filecoin'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range('D35') '5.347'
$ws.Range('E35').Value = '  +2.07%  '

# Row 36
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range('D36') '1.372'
$ws.Range('E36').Value = '  +4.71%  '

# Row 37
Set-TextValue $ws.Range('D37') '0.06149'
$ws.Range('E37').Value = '  +3.69%  '

# Row 38
Set-TextValue $ws.Range('D38') '0.02251'
$ws.Range('E38').Value = '  +3.13%  '

# Row 39
Set-TextValue $ws.Range('D39') '8.273'
$ws.Range('E39').Value = '  +3.80%  '

# Row 40
Set-TextValue $ws.Range('D40') '1.178'
$ws.Range('E40').Value = '  +1.94%  '

# Row 41
Set-TextValue $ws.Range('D41') '0.5910'
$ws.Range('E41').Value = '  +2.27%  '

# Row 42
$ws.Range('E42').Value = '  -0.53%  '

# Row 43
$ws.Range('E43').Value = '  +2.13%  '

# Row 44
Set-TextValue $ws.Range('D44') '10.27'

# Row 45
Set-TextValue $ws.Range('D45') '1.273'
$ws.Range('E45').Value = '  -0.26%  '

# Row 46
Set-TextValue $ws.Range('D46') '0.07621'
$ws.Range('E46').Value = '  +14.72%  '

# Row 47
Set-TextValue $ws.Range('D47') '0.5533'
$ws.Range('E47').Value = '  +0.98%  '

# Row 48
Set-TextValue $ws.Range('D48') '12.05'
$ws.Range('E48').Value = '  +1.00%  '

# Row 49
Set-TextValue $ws.Range('D49') '1.938'
$ws.Range('E49').Value = '  +3.99%  '

# Row 50
Set-TextValue $ws.Range('D50') '2.063'
$ws.Range('E50').Value = '  +13.50%  '

# Row 51
Set-TextValue $ws.Range('D51') '111.83'
$ws.Range('E51').Value = '  +1.81%  '
